$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 645 ("「車のバッテリーが要る！」" post), shifting all
# subsequent rows up by one.
$ws.Rows.Item(645).Delete()
